# The document consists of a series of list items (each with a label
# paragraph and a following paragraph holding a screenshot image) that
# illustrate several ML models, followed by one trailing empty
# "ListParagraph"-styled paragraph before the section break.
#
# This edit removes everything except that final trailing empty
# paragraph: all of the leading picture, the "Logistic regression" /
# "KNN" / "Naive Bayes" / "Decision Tree" bullets and their associated
# screenshots/spacer paragraphs are deleted, leaving just the last
# empty paragraph before the section properties.

$d = $word.ActiveDocument

$paragraphs = $d.Paragraphs
$count = $paragraphs.Count

if ($count -gt 1) {
    $lastPara = $paragraphs.Item($count)
    $deleteRange = $d.Range(0, $lastPara.Range.Start)
    $deleteRange.Delete()
}
